$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving its text data-type and original style,
# even when the new value looks like a plain number (e.g. "34.26") which Excel
# would otherwise auto-convert to a numeric cell.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "47.770.98"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.477.00"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "315.96"
$ws.Range("E5").Value = "  -2.15%  "
Set-TextValue $ws.Range("D6") "103.65"
$ws.Range("E6").Value = "  -5.51%  "
$ws.Range("E7").Value = "  -3.27%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue $ws.Range("D9") "0.532"
$ws.Range("E9").Value = "  -4.15%  "
Set-TextValue $ws.Range("D10") "38.53"
$ws.Range("E10").Value = "  -5.10%  "
Set-TextValue $ws.Range("D11") "20.40"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("E12").Value = "  -3.52%  "
Set-TextValue $ws.Range("D14") "6.99"
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("D15").Value = "2.865.92"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "2.495.64"
$ws.Range("E16").Value = "  -1.15%  "
Set-TextValue $ws.Range("D17") "0.819"
$ws.Range("E17").Value = "  -4.25%  "
$ws.Range("D18").Value = "47.712.35"
$ws.Range("E18").Value = "  -0.94%  "
Set-TextValue $ws.Range("D19") "2.90"
$ws.Range("E19").Value = "  +7.51%  "
Set-TextValue $ws.Range("D20") "12.58"
$ws.Range("E20").Value = "  -6.56%  "
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("E22").Value = "  -2.81%  "
Set-TextValue $ws.Range("D23") "279.51"
$ws.Range("E23").Value = "  +5.39%  "
Set-TextValue $ws.Range("D24") "70.67"
$ws.Range("E24").Value = "  -1.94%  "
Set-TextValue $ws.Range("D25") "2.47"
$ws.Range("E25").Value = "  -3.85%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  +0.14%  "
Set-TextValue $ws.Range("D27") "25.59"
$ws.Range("E27").Value = "  -2.02%  "
Set-TextValue $ws.Range("D28") "2.21"
$ws.Range("E28").Value = "  +0.21%  "
Set-TextValue $ws.Range("D29") "9.51"
$ws.Range("E29").Value = "  -6.05%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D30") "34.26"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D31") "0.135"
$ws.Range("E31").Value = "  -5.90%  "
Set-TextValue $ws.Range("D32") "49.15"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  -0.06%  "
Set-TextValue $ws.Range("D34") "18.80"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("E35").Value = "  -3.73%  "
Set-TextValue $ws.Range("D36") "0.0765"
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("E37").Value = "  -4.02%  "
Set-TextValue $ws.Range("D38") "4.45"
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("E39").Value = "  -5.59%  "
Set-TextValue $ws.Range("D40") "122.80"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "21.99"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D43") "2.19"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "1.988.14"
$ws.Range("E45").Value = "  -1.59%  "
Set-TextValue $ws.Range("D46") "3.10"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -2.06%  "
Set-TextValue $ws.Range("D48") "1.96"
$ws.Range("E48").Value = "  -4.11%  "
$ws.Range("E49").Value = "  -3.02%  "
Set-TextValue $ws.Range("D50") "5.06"
$ws.Range("E50").Value = "  -3.41%  "
Set-TextValue $ws.Range("D51") "78.79"
$ws.Range("E51").Value = "  -0.64%  "
